# Auto-generated edit script: updates per-row H..N values to match the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 805.5484
$ws.Range("I33").Value = 190.625
$ws.Range("J33").Value = 2913.8572
$ws.Range("K33").Value = 190.625
$ws.Range("L33").Value = 2913.8572
$ws.Range("M33").Value = 38.375
$ws.Range("N33").Value = -3371.8572
$ws.Range("H96").Value = 774.5333000000001
$ws.Range("I96").Value = 595.1818
$ws.Range("J96").Value = 1267.75
$ws.Range("K96").Value = 1785.5454
$ws.Range("L96").Value = 3803.25
$ws.Range("M96").Value = -412.5454
$ws.Range("N96").Value = -6549.25
$ws.Range("H132").Value = 8071623.5
$ws.Range("I132").Value = 8340377.5
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 25021132.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -25018602.5
$ws.Range("N132").Value = -32060
$ws.Range("H137").Value = 1471.0233
$ws.Range("I137").Value = 1628.8
$ws.Range("J137").Value = 1423.2122
$ws.Range("K137").Value = 4886.4
$ws.Range("L137").Value = 4269.6366
$ws.Range("M137").Value = -2336.4
$ws.Range("N137").Value = -9369.6366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 19514.27
$ws.Range("I32").Value = 3104.75
$ws.Range("K32").Value = 3104.75
$ws.Range("M32").Value = -2817.75
$ws.Range("H45").Value = 2033.4138
$ws.Range("I45").Value = 1789.2778
$ws.Range("J45").Value = 2432.9092
$ws.Range("K45").Value = 1789.2778
$ws.Range("L45").Value = 2432.9092
$ws.Range("M45").Value = -1412.2778
$ws.Range("N45").Value = -3186.9092
$ws.Range("H61").Value = 1894.5853
$ws.Range("I61").Value = 1041.5
$ws.Range("J61").Value = 2440.56
$ws.Range("K61").Value = 1041.5
$ws.Range("L61").Value = 2440.56
$ws.Range("M61").Value = -829.5
$ws.Range("N61").Value = -2864.56
$ws.Range("H74").Value = 1785.683
$ws.Range("I74").Value = 1038.0476
$ws.Range("J74").Value = 2570.7
$ws.Range("K74").Value = 1038.0476
$ws.Range("L74").Value = 2570.7
$ws.Range("M74").Value = -164.0476000000001
$ws.Range("N74").Value = -4318.7
$ws.Range("H77").Value = 1785.683
$ws.Range("I77").Value = 1038.0476
$ws.Range("J77").Value = 2570.7
$ws.Range("K77").Value = 5190.238
$ws.Range("L77").Value = 12853.5
$ws.Range("M77").Value = -822.2380000000003
$ws.Range("N77").Value = -21589.5
$ws.Range("H136").Value = 1894.5853
$ws.Range("I136").Value = 1041.5
$ws.Range("J136").Value = 2440.56
$ws.Range("K136").Value = 3124.5
$ws.Range("L136").Value = 7321.68
$ws.Range("M136").Value = -574.5
$ws.Range("N136").Value = -12421.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19383.5
$ws.Range("J35").Value = 19383.5
$ws.Range("L35").Value = 19383.5
$ws.Range("N35").Value = -20003.5
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1117
$ws.Range("I16").Value = 1051
$ws.Range("J16").Value = 1150
$ws.Range("K16").Value = 1051
$ws.Range("L16").Value = 1150
$ws.Range("M16").Value = -764
$ws.Range("N16").Value = -1724
$ws.Range("H31").Value = 16802.135
$ws.Range("I31").Value = 39438.5
$ws.Range("J31").Value = 2447.366
$ws.Range("K31").Value = 39438.5
$ws.Range("L31").Value = 2447.366
$ws.Range("M31").Value = -39143.5
$ws.Range("N31").Value = -3037.366
$ws.Range("H34").Value = 16802.135
$ws.Range("I34").Value = 39438.5
$ws.Range("J34").Value = 2447.366
$ws.Range("K34").Value = 39438.5
$ws.Range("L34").Value = 2447.366
$ws.Range("M34").Value = -39236.5
$ws.Range("N34").Value = -2851.366
$ws.Range("H107").Value = 613.1667
$ws.Range("I107").Value = 558.0476
$ws.Range("J107").Value = 741.7778
$ws.Range("K107").Value = 558.0476
$ws.Range("L107").Value = 741.7778
$ws.Range("M107").Value = 1361.9524
$ws.Range("N107").Value = -4581.7778
$ws.Range("H113").Value = 1117
$ws.Range("I113").Value = 1051
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 1051
$ws.Range("L113").Value = 1150
$ws.Range("M113").Value = 1119
$ws.Range("N113").Value = -5490
$ws.Range("H132").Value = 2110.2812
$ws.Range("I132").Value = 2077.6333
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 6232.8999
$ws.Range("L132").Value = 7800
$ws.Range("M132").Value = -3702.8999
$ws.Range("N132").Value = -12860
$ws.Range("H134").Value = 1497.3704
$ws.Range("I134").Value = 1225
$ws.Range("K134").Value = 3675
$ws.Range("M134").Value = -1140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1981.9863
$ws.Range("I68").Value = 1281.4688
$ws.Range("J68").Value = 2528.7317
$ws.Range("K68").Value = 3844.4064
$ws.Range("L68").Value = 7586.195099999999
$ws.Range("M68").Value = -3033.4064
$ws.Range("N68").Value = -9208.195099999999
$ws.Range("H71").Value = 1981.9863
$ws.Range("I71").Value = 1281.4688
$ws.Range("J71").Value = 2528.7317
$ws.Range("K71").Value = 11533.2192
$ws.Range("L71").Value = 22758.5853
$ws.Range("M71").Value = -7477.219200000001
$ws.Range("N71").Value = -30870.5853
$ws.Range("H131").Value = 1328.7627
$ws.Range("I131").Value = 935.55554
$ws.Range("J131").Value = 1399.54
$ws.Range("K131").Value = 2806.66662
$ws.Range("L131").Value = 4198.62
$ws.Range("M131").Value = 2233.33338
$ws.Range("N131").Value = -14278.62

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 44145
$ws.Range("J24").Value = 44145
$ws.Range("L24").Value = 44145
$ws.Range("N24").Value = -44491
$ws.Range("H126").Value = 4904376
$ws.Range("I126").Value = 3198.6667
$ws.Range("J126").Value = 9805553
$ws.Range("K126").Value = 9596.000100000001
$ws.Range("L126").Value = 29416659
$ws.Range("M126").Value = -7126.000100000001
$ws.Range("N126").Value = -29421599
$ws.Range("H132").Value = 2594.1177
$ws.Range("I132").Value = 2291.6667
$ws.Range("J132").Value = 3320
$ws.Range("K132").Value = 6875.000100000001
$ws.Range("L132").Value = 9960
$ws.Range("M132").Value = -4345.000100000001
$ws.Range("N132").Value = -15020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3964.0476
$ws.Range("I132").Value = 4327.625
$ws.Range("J132").Value = 2800.6
$ws.Range("K132").Value = 12982.875
$ws.Range("L132").Value = 8401.799999999999
$ws.Range("M132").Value = -10452.875
$ws.Range("N132").Value = -13461.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 50010
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 50010
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 50010
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -50596
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H123").Value = 37330
$ws.Range("J123").Value = 37330
$ws.Range("L123").Value = 37330
$ws.Range("N123").Value = -47130
$ws.Range("H132").Value = 1918.283
$ws.Range("I132").Value = 2158.5527
$ws.Range("J132").Value = 1309.6
$ws.Range("K132").Value = 6475.658100000001
$ws.Range("L132").Value = 3928.8
$ws.Range("M132").Value = -3945.658100000001
$ws.Range("N132").Value = -8988.799999999999
$ws.Range("H136").Value = 962.4138
$ws.Range("I136").Value = 495.95456
$ws.Range("J136").Value = 2428.4285
$ws.Range("K136").Value = 1487.86368
$ws.Range("L136").Value = 7285.2855
$ws.Range("M136").Value = 1062.13632
$ws.Range("N136").Value = -12385.2855

